$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Set the default value for the "dataImporterConfigurationFile" property (row 12, column B)
$ws.Range("B12").Value = "esqlabs_dataImporter_configuration.xml"

# Update selection to the newly edited cell
$ws.Range("B12").Select()

# Update workbook window position to match the author's saved view
$excel.ActiveWindow.Left = 7575
$excel.ActiveWindow.Top = 2475
